$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns C (runs), D (balls), E (fours), F (sixes).
# Derived from the authoritative diff; only the cells whose value actually
# changes are touched here - everything else in the sheet is left as-is.
$updates = @{
    2  = @{ C = "11"; D = "13"; E = "1" }
    3  = @{ C = "88"; D = "38"; E = "7"; F = "6" }
    4  = @{ C = "65"; D = "50"; E = "6"; F = "2" }
    7  = @{ C = "25"; D = "29"; E = "1"; F = "1" }
    8  = @{ C = "14"; D = "12"; E = "0"; F = "1" }
    9  = @{ C = "7";  D = "9";  E = "0"; F = "0" }
    11 = @{ C = "26"; D = "22"; E = "1"; F = "0" }
    12 = @{ C = "17"; D = "21"; E = "2"; F = "0" }
    13 = @{ C = "39"; D = "32"; E = "0"; F = "3" }
    14 = @{ C = "23"; D = "23"; E = "1"; F = "1" }
    15 = @{ C = "53"; D = "43"; E = "3"; F = "2" }
    16 = @{ C = "7";  D = "12"; F = "0" }
    17 = @{ C = "47"; D = "38"; E = "5" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        # Prefix with an apostrophe so Excel stores the numeric-looking
        # value as text (quote-prefixed), matching the original workbook's
        # string-typed "runs"/"balls"/"fours"/"sixes" cells instead of
        # converting them to numbers.
        $ws.Range("$col$row").Value = "'" + $cols[$col]
    }
}
